# Auto-generated edit script updating cryptos list values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.689.11'
$ws.Range('E2').Value = '  +1.38%  '
$ws.Range('D3').Value = '3.607.24'
$ws.Range('E3').Value = '  +2.20%  '
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  +0.21%  '
$ws.Range('D5').Value = '''202.04'
$ws.Range('E5').Value = '  +4.15%  '
$ws.Range('D6').Value = '''597.82'
$ws.Range('E6').Value = '  -1.51%  '
$ws.Range('D7').Value = '''0.628'
$ws.Range('E7').Value = '  +0.86%  '
$ws.Range('D8').Value = '''1.00'
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('D9').Value = '''0.215'
$ws.Range('E9').Value = '  +7.02%  '
$ws.Range('D10').Value = '''0.645'
$ws.Range('E10').Value = '  +0.01%  '
$ws.Range('D11').Value = '''53.97'
$ws.Range('E11').Value = '  +1.40%  '
$ws.Range('D12').Value = '''0.0000302'
$ws.Range('E12').Value = '  +0.20%  '
$ws.Range('D13').Value = '''9.66'
$ws.Range('E13').Value = '  +2.08%  '
$ws.Range('D14').Value = '4.180.16'
$ws.Range('E14').Value = '  +2.30%  '
$ws.Range('D15').Value = '''680.83'
$ws.Range('E15').Value = '  +14.92%  '
$ws.Range('D16').Value = '70.771.21'
$ws.Range('E16').Value = '  +1.36%  '
$ws.Range('B17').Value = 'Chainlink'
$ws.Range('C17').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D17').Value = '''19.20'
$ws.Range('E17').Value = '  +1.12%  '
$ws.Range('B18').Value = 'Uniswap'
$ws.Range('C18').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D18').Value = '''12.80'
$ws.Range('E18').Value = '  +0.67%  '
$ws.Range('B19').Value = 'WrappedEther'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D19').Value = '3.602.83'
$ws.Range('E19').Value = '  +1.35%  '
$ws.Range('E20').Value = '  +0.60%  '
$ws.Range('E21').Value = '  +1.98%  '
$ws.Range('D22').Value = '''18.76'
$ws.Range('E22').Value = '  +5.56%  '
$ws.Range('D23').Value = '''110.47'
$ws.Range('E23').Value = '  +7.93%  '
$ws.Range('D24').Value = '''5.28'
$ws.Range('E24').Value = '  +3.15%  '
$ws.Range('D25').Value = '''4.59'
$ws.Range('E25').Value = '  -1.02%  '
$ws.Range('E26').Value = '  -0.15%  '
$ws.Range('D27').Value = '''10.65'
$ws.Range('E27').Value = '  -0.97%  '
$ws.Range('D28').Value = '''6.01'
$ws.Range('E28').Value = '  -0.49%  '
$ws.Range('D29').Value = '''10.15'
$ws.Range('E29').Value = '  +6.59%  '
$ws.Range('D30').Value = '''34.46'
$ws.Range('E30').Value = '  +4.03%  '
$ws.Range('D31').Value = '''4.49'
$ws.Range('E31').Value = '  +6.78%  '
$ws.Range('D32').Value = '''7.19'
$ws.Range('E32').Value = '  +2.32%  '
$ws.Range('D33').Value = '''12.32'
$ws.Range('E33').Value = '  +0.10%  '
$ws.Range('E34').Value = '  +0.44%  '
$ws.Range('D35').Value = '''63.54'
$ws.Range('E35').Value = '  +0.38%  '
$ws.Range('D36').Value = '3.888.76'
$ws.Range('E36').Value = '  +2.17%  '
$ws.Range('D37').Value = '0.0₃0851'
$ws.Range('E37').Value = '  +5.20%  '
$ws.Range('D38').Value = '''0.999'
$ws.Range('E38').Value = '  -0.15%  '
$ws.Range('D39').Value = '''515.80'
$ws.Range('E39').Value = '  +1.05%  '
$ws.Range('E40').Value = '  -4.98%  '
$ws.Range('E41').Value = '  +0.59%  '
$ws.Range('D42').Value = '''36.87'
$ws.Range('E42').Value = '  +0.90%  '
$ws.Range('D43').Value = '''0.386'
$ws.Range('E43').Value = '  -1.39%  '
$ws.Range('E44').Value = '  +3.40%  '
$ws.Range('D45').Value = '''0.0467'
$ws.Range('E45').Value = '  +4.53%  '
$ws.Range('D46').Value = '''3.08'
$ws.Range('E46').Value = '  +9.78%  '
$ws.Range('E47').Value = '  +3.40%  '
$ws.Range('E48').Value = '  +1.71%  '
$ws.Range('E49').Value = '  +2.10%  '
$ws.Range('E50').Value = '  -0.18%  '
$ws.Range('D51').Value = '''1.81'
$ws.Range('E51').Value = '  +23.07%  '
